$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.081742303013206
$ws.Range("D2").Value = 1.076586557308225
$ws.Range("E2").Value = 1.095574217330284
$ws.Range("F2").Value = 1.102888618657688
$ws.Range("I2").Value = 1.057874511200575
$ws.Range("J2").Value = 1.086615379532172
$ws.Range("K2").Value = 1.079270098372152
$ws.Range("L2").Value = 1.098208676608113
$ws.Range("M2").Value = 1.105504656545144
$ws.Range("N2").Value = 1.088158498007634

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.083361597210684
$ws.Range("D3").Value = 1.077872436652297
$ws.Range("E3").Value = 1.09722382004255
$ws.Range("F3").Value = 1.104680696191943
$ws.Range("I3").Value = 1.058436170805519
$ws.Range("J3").Value = 1.087893416690349
$ws.Range("K3").Value = 1.080372365266487
$ws.Range("L3").Value = 1.099677380441298
$ws.Range("M3").Value = 1.107116854112412
$ws.Range("N3").Value = 1.089438350125169

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.084407060054298
$ws.Range("D4").Value = 1.07870226403431
$ws.Range("E4").Value = 1.098289236394857
$ws.Range("F4").Value = 1.105838510232278
$ws.Range("I4").Value = 1.058796994135388
$ws.Range("J4").Value = 1.088717647953479
$ws.Range("K4").Value = 1.081082828493557
$ws.Range("L4").Value = 1.100625253647451
$ws.Range("M4").Value = 1.108157808520277
$ws.Range("N4").Value = 1.090263751891233

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.084846026360484
$ws.Range("D5").Value = 1.079050599607273
$ws.Range("E5").Value = 1.098736672682675
$ws.Range("F5").Value = 1.106324841039538
$ws.Range("I5").Value = 1.058948064319411
$ws.Range("J5").Value = 1.089063506825489
$ws.Range("K5").Value = 1.081380850141989
$ws.Range("L5").Value = 1.101023156819277
$ws.Range("M5").Value = 1.108594899235992
$ws.Range("N5").Value = 1.090610101922511

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.084919699026273
$ws.Range("D6").Value = 1.079109056201136
$ws.Range("E6").Value = 1.098811772267194
$ws.Range("F6").Value = 1.106406474176779
$ws.Range("I6").Value = 1.058973393441871
$ws.Range("J6").Value = 1.089121540316782
$ws.Range("K6").Value = 1.081430850980446
$ws.Range("L6").Value = 1.10108993263436
$ws.Range("M6").Value = 1.108668258133202
$ws.Range("N6").Value = 1.090668217828018

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.084412927674502
$ws.Range("D7").Value = 1.078706920558256
$ws.Range("E7").Value = 1.098295216871448
$ws.Range("F7").Value = 1.105845010219442
$ws.Range("I7").Value = 1.058799015172382
$ws.Range("J7").Value = 1.088722271870148
$ws.Range("K7").Value = 1.081086813244547
$ws.Range("L7").Value = 1.10063057272479
$ws.Range("M7").Value = 1.108163651001455
$ws.Range("N7").Value = 1.090268382374393

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.082290040047845
$ws.Range("D8").Value = 1.077021591411214
$ws.Range("E8").Value = 1.09613212629011
$ws.Range("F8").Value = 1.103494635395549
$ws.Range("I8").Value = 1.058064869362905
$ws.Range("J8").Value = 1.087047871837911
$ws.Range("K8").Value = 1.079643194350598
$ws.Range("L8").Value = 1.098705550813404
$ws.Range("M8").Value = 1.106049977789222
$ws.Range("N8").Value = 1.088591604502055

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.078530899259348
$ws.Range("D9").Value = 1.074034433814867
$ws.Range("E9").Value = 1.092304775387008
$ws.Range("F9").Value = 1.099338806442413
$ws.Range("I9").Value = 1.056751037628274
$ws.Range("J9").Value = 1.084075941121318
$ws.Range("K9").Value = 1.077077733932063
$ws.Range("L9").Value = 1.095293994276627
$ws.Range("M9").Value = 1.102307725075822
$ws.Range("N9").Value = 1.085615453302959

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.076011763500548
$ws.Range("D10").Value = 1.07203078110142
$ws.Range("E10").Value = 1.089741956141397
$ws.Range("F10").Value = 1.096557975850309
$ws.Range("I10").Value = 1.055861306706973
$ws.Range("J10").Value = 1.082079678872972
$ws.Range("K10").Value = 1.075352399600116
$ws.Range("L10").Value = 1.093005924872616
$ws.Range("M10").Value = 1.099800301336341
$ws.Range("N10").Value = 1.083616356133245

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.074917697759818
$ws.Range("D11").Value = 1.071160163570978
$ws.Range("E11").Value = 1.08862940391355
$ws.Range("F11").Value = 1.095351239582842
$ws.Range("I11").Value = 1.055472697184122
$ws.Range("J11").Value = 1.081211595575708
$ws.Range("K11").Value = 1.07460163631708
$ws.Range("L11").Value = 1.092011773250213
$ws.Range("M11").Value = 1.098711417991139
$ws.Range("N11").Value = 1.082747040058128

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.07451080898668
$ws.Range("D12").Value = 1.070836313182338
$ws.Range("E12").Value = 1.088215713230116
$ws.Range("F12").Value = 1.094902596073327
$ws.Range("I12").Value = 1.055327841494909
$ws.Range("J12").Value = 1.080888585517637
$ws.Range("K12").Value = 1.074322206622055
$ws.Range("L12").Value = 1.091641977459794
$ws.Range("M12").Value = 1.098306471164771
$ws.Range("N12").Value = 1.082423571288726

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.074598111041184
$ws.Range("D13").Value = 1.070905801363434
$ws.Range("E13").Value = 1.088304471379279
$ws.Range("F13").Value = 1.094998850324082
$ws.Range("I13").Value = 1.055358936645323
$ws.Range("J13").Value = 1.080957898030671
$ws.Range("K13").Value = 1.07438217086225
$ws.Range("L13").Value = 1.091721323807599
$ws.Range("M13").Value = 1.098393355892201
$ws.Range("N13").Value = 1.082492982233479

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.07488407456934
$ws.Range("D14").Value = 1.071133403526945
$ws.Range("E14").Value = 1.088595217151061
$ws.Range("F14").Value = 1.095314162976786
$ws.Range("I14").Value = 1.055460733779152
$ws.Range("J14").Value = 1.081184907059506
$ws.Range("K14").Value = 1.074578550107733
$ws.Range("L14").Value = 1.091981216568709
$ws.Range("M14").Value = 1.098677954974597
$ws.Range("N14").Value = 1.082720313641172

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.075060198932287
$ws.Range("D15").Value = 1.071273574905575
$ws.Range("E15").Value = 1.088774296609384
$ws.Range("F15").Value = 1.095508383015994
$ws.Range("I15").Value = 1.055523386765792
$ws.Range("J15").Value = 1.081324699576599
$ws.Range("K15").Value = 1.074699470979491
$ws.Range("L15").Value = 1.092141275318606
$ws.Range("M15").Value = 1.098853240915632
$ws.Range("N15").Value = 1.082860304679675

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.076084303113965
$ws.Range("D16").Value = 1.072088496571133
$ws.Range("E16").Value = 1.089815731707665
$ws.Range("F16").Value = 1.096638006558172
$ws.Range("I16").Value = 1.05588702639009
$ws.Range("J16").Value = 1.082137212043187
$ws.Range("K16").Value = 1.07540214689235
$ws.Range("L16").Value = 1.093071830697212
$ws.Range("M16").Value = 1.099872499375354
$ws.Range("N16").Value = 1.083673971007161

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.076725813247148
$ws.Range("D17").Value = 1.072598859010022
$ws.Range("E17").Value = 1.090468228122492
$ws.Range("F17").Value = 1.097345878337911
$ws.Range("I17").Value = 1.056114227129728
$ws.Range("J17").Value = 1.082645884306705
$ws.Range("K17").Value = 1.075841924379324
$ws.Range("L17").Value = 1.093654623824943
$ws.Range("M17").Value = 1.100510999844733
$ws.Range("N17").Value = 1.084183365643639

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.077099681326643
$ws.Range("D18").Value = 1.072896253831976
$ws.Range("E18").Value = 1.090848545656751
$ws.Range("F18").Value = 1.097758516336698
$ws.Range("I18").Value = 1.056246426608638
$ws.Range("J18").Value = 1.082942229094223
$ws.Range("K18").Value = 1.076098084527285
$ws.Range("L18").Value = 1.093994229753097
$ws.Range("M18").Value = 1.100883123129409
$ws.Range("N18").Value = 1.084480131274747

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.07722710784689
$ws.Range("D19").Value = 1.072997608735409
$ws.Range("E19").Value = 1.09097817827537
$ws.Range("F19").Value = 1.097899172987669
$ws.Range("I19").Value = 1.056291448615529
$ws.Range("J19").Value = 1.083043215105183
$ws.Range("K19").Value = 1.076185368656727
$ws.Range("L19").Value = 1.094109971432475
$ws.Range("M19").Value = 1.101009956524195
$ws.Range("N19").Value = 1.084581260697425

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.076657017851565
$ws.Range("D20").Value = 1.072544132118481
$ws.Range("E20").Value = 1.090398249723877
$ws.Range("F20").Value = 1.097269956528157
$ws.Range("I20").Value = 1.056089884076053
$ws.Range("J20").Value = 1.082591345400695
$ws.Range("K20").Value = 1.075794777166272
$ws.Range("L20").Value = 1.093592129596683
$ws.Range("M20").Value = 1.100442526249624
$ws.Range("N20").Value = 1.084128749286127

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.074799879463025
$ws.Range("D21").Value = 1.071066393260863
$ws.Range("E21").Value = 1.088509612001499
$ws.Range("F21").Value = 1.095221322638043
$ws.Range("I21").Value = 1.055430771150424
$ws.Range("J21").Value = 1.081118074252499
$ws.Range("K21").Value = 1.074520736927603
$ws.Range("L21").Value = 1.091904699182932
$ws.Range("M21").Value = 1.098594161192589
$ws.Range("N21").Value = 1.082653385923914

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.073629299692934
$ws.Range("D22").Value = 1.070134588944907
$ws.Range("E22").Value = 1.087319602418754
$ws.Range("F22").Value = 1.093930897729765
$ws.Range("I22").Value = 1.055013414122118
$ws.Range("J22").Value = 1.080188494150607
$ws.Range("K22").Value = 1.073716436351838
$ws.Range("L22").Value = 1.090840708998611
$ws.Range("M22").Value = 1.097429195920721
$ws.Range("N22").Value = 1.081722485711653

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.074250127887802
$ws.Range("D23").Value = 1.070628814556449
$ws.Range("E23").Value = 1.087950695045402
$ws.Range("F23").Value = 1.094615205696292
$ws.Range("I23").Value = 1.055234944125639
$ws.Range("J23").Value = 1.080681596418087
$ws.Range("K23").Value = 1.074143123625961
$ws.Range("L23").Value = 1.09140504225758
$ws.Range("M23").Value = 1.098047038341399
$ws.Range("N23").Value = 1.082216288240912

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.076688104499565
$ws.Range("D24").Value = 1.072568861744167
$ws.Range("E24").Value = 1.09042987079231
$ws.Range("F24").Value = 1.097304263103498
$ws.Range("I24").Value = 1.056100884651836
$ws.Range("J24").Value = 1.082615990282081
$ws.Range("K24").Value = 1.075816082058198
$ws.Range("L24").Value = 1.093620369060184
$ws.Range("M24").Value = 1.10047346745682
$ws.Range("N24").Value = 1.084153429166072

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.07950497354232
$ws.Range("D25").Value = 1.074808800037769
$ws.Range("E25").Value = 1.093296169600318
$ws.Range("F25").Value = 1.100414943864302
$ws.Range("I25").Value = 1.057093113882934
$ws.Range("J25").Value = 1.084846853742616
$ws.Range("K25").Value = 1.077743579297584
$ws.Range("L25").Value = 1.096178330567218
$ws.Range("M25").Value = 1.103277354635422
$ws.Range("N25").Value = 1.086387460708604
